# Auto-generated script applying cryptos.xlsx data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.375.87'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -2.40%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.992.40'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -5.93%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.007'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '330.75'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -4.60%  '
$ws.Range('E6').Value = '  +0.08%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4945'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -4.89%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4198'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -5.86%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '52.10'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.86%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.08845'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -5.54%  '
$ws.Range('E11').Value = '  -5.46%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '23.26'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -8.43%  '
$ws.Range('B13').Value = 'Chainlink'
$ws.Range('C13').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '8.032'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -7.31%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.969.95'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -7.14%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.503'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -6.51%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '96.27'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -6.13%  '
$ws.Range('E17').Value = '  +0.05%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001104'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -5.42%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06633'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.96%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '19.66'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -8.83%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.007'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.11%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.961'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -5.22%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '29.409.49'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.39%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.84'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -7.16%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.290'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.34%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '157.60'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.96%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.582'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.08%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.52'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -7.44%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.344'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -7.62%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '127.33'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -5.03%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.054'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -8.63%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09936'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -6.00%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.564'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -12.47%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.842'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -6.80%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.774'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.93%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '9.598'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -11.00%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02452'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.06355'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -7.43%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.286'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.34%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '11.75'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -7.48%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6494'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -8.63%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.2066'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -8.07%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.006'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.09%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6323'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -7.84%  '
$ws.Range('B45').Value = 'NEARProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.205'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -7.60%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '13.39'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -8.83%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.276'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.49%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.534'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.67%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.00000000334'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.41%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06988'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.69%  '
$ws.Range('B51').Value = 'ThetaToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.145'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.32%  '
